# Applies the "Upload new version with timestamp" edit:
#  - A new item row "FERROTRON 30 CAPS" (balance 0:0, price 138.00, sell price
#    138.0000, transactions 1:0) is inserted into the shortage list right
#    after "EZACARD ..." (row 13), pushing the following rows
#    ("OXITROPIL ...", "SELGON ...") down by one.
#  - The old "SIDERAL FOLIC 20 ORODISPERSIBLE STICKS" row is removed from
#    the list (its slot is effectively consumed by the downward shift).
#  - The printed total in P23 is updated to reflect the new price total.
#  - The generated timestamp text in A24 is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $txt) {
    # Some columns are formatted as real numbers (e.g. numFmtId "0.00").
    # Assigning a numeric-looking string to those cells makes Excel coerce
    # it into a genuine number, which loses the original textual
    # representation (e.g. "138.0000" -> 138). To keep the cell's original
    # literal text (and its original style/number format), temporarily
    # switch to a text format, assign the value, then restore the format.
    $orig = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $txt
    $rng.NumberFormat = $orig
}

# Row 13: was OXITROPIL 1200 MG 60 TAB -> now FERROTRON 30 CAPS
Set-TextValue $ws.Range("C13") "FERROTRON 30 CAPS"
Set-TextValue $ws.Range("H13") "0:0"
Set-TextValue $ws.Range("N13") "138.00"
Set-TextValue $ws.Range("P13") "138.0000"
Set-TextValue $ws.Range("Q13") "1:0"

# Row 14: was SELGON 10MG 6 INFANT SUPP. -> now OXITROPIL 1200 MG 60 TAB
Set-TextValue $ws.Range("C14") "OXITROPIL 1200 MG 60 TAB"
Set-TextValue $ws.Range("H14") "1:1"
Set-TextValue $ws.Range("N14") "123.00"
Set-TextValue $ws.Range("P14") "19.6800"
Set-TextValue $ws.Range("Q14") "0:1"

# Row 15: was SIDERAL FOLIC 20 ORODISPERSIBLE STICKS -> now SELGON 10MG 6 INFANT SUPP.
Set-TextValue $ws.Range("C15") "SELGON 10MG 6 INFANT SUPP."
Set-TextValue $ws.Range("H15") "3:0"
Set-TextValue $ws.Range("N15") "15.00"
Set-TextValue $ws.Range("P15") "15.0000"
Set-TextValue $ws.Range("Q15") "1:0"

# Updated printed total (975.215 - 180.00 (removed) + 138.00 (added) = 933.215)
$ws.Range("P23").Value = 933.21500000000003

# Updated generation timestamp
$ws.Range("A24").Value = "Monday, 8 September, 2025 11:06 AM"

Write-Host "Applied FERROTRON insertion / SIDERAL FOLIC removal and timestamp update"
